# Updated symbol list on Wed Jan 18 15:36:15 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the crypto symbol table.
# Values are stored as plain text in the sheet (not numbers/percentages), so we
# temporarily force a Text number format before writing each value and then
# clear the formatting again to leave styling untouched (matching the source
# workbook, where these cells carry no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "301.58" }
    @{ Cell = "E2"; Value = "-0.58%" }
    @{ Cell = "D3"; Value = "32.65" }
    @{ Cell = "E3"; Value = "1.79%" }
    @{ Cell = "D4"; Value = "5.080" }
    @{ Cell = "E4"; Value = "-0.86%" }
    @{ Cell = "D5"; Value = "0.07725" }
    @{ Cell = "E5"; Value = "-1.57%" }
    @{ Cell = "D6"; Value = "2.038" }
    @{ Cell = "E6"; Value = "-9.85%" }
    @{ Cell = "D7"; Value = "7.915" }
    @{ Cell = "E7"; Value = "0.94%" }
    @{ Cell = "D8"; Value = "3.802" }
    @{ Cell = "E8"; Value = "-0.04%" }
    @{ Cell = "D9"; Value = "0.9257" }
    @{ Cell = "E9"; Value = "0.01%" }
    @{ Cell = "D10"; Value = "0.1763" }
    @{ Cell = "E10"; Value = "-0.30%" }
    @{ Cell = "D11"; Value = "0.08157" }
    @{ Cell = "E11"; Value = "4.60%" }
    @{ Cell = "D12"; Value = "0.08592" }
    @{ Cell = "E12"; Value = "-3.52%" }
    @{ Cell = "D13"; Value = "0.03052" }
    @{ Cell = "E13"; Value = "-1.12%" }
    @{ Cell = "D14"; Value = "0.09977" }
    @{ Cell = "E14"; Value = "-0.68%" }
    @{ Cell = "D15"; Value = "0.001523" }
    @{ Cell = "E15"; Value = "0.55%" }
    @{ Cell = "D16"; Value = "0.005899" }
    @{ Cell = "E16"; Value = "0.36%" }
    @{ Cell = "D18"; Value = "3.476" }
    @{ Cell = "E18"; Value = "0.49%" }
    @{ Cell = "E19"; Value = "-4.11%" }
    @{ Cell = "D20"; Value = "0.3334" }
    @{ Cell = "E20"; Value = "1.85%" }
    @{ Cell = "E21"; Value = "-0.13%" }
    @{ Cell = "D22"; Value = "4.409" }
    @{ Cell = "E22"; Value = "3.29%" }
    @{ Cell = "D23"; Value = "0.1977" }
    @{ Cell = "E23"; Value = "10.15%" }
    @{ Cell = "D24"; Value = "0.04545" }
    @{ Cell = "E24"; Value = "-1.40%" }
    @{ Cell = "D25"; Value = "0.001231" }
    @{ Cell = "E25"; Value = "-1.73%" }
    @{ Cell = "D26"; Value = "0.004153" }
    @{ Cell = "E26"; Value = "-7.70%" }
    @{ Cell = "D27"; Value = "0.0001252" }
    @{ Cell = "E27"; Value = "0.12%" }
    @{ Cell = "D39"; Value = "0.01735" }
    @{ Cell = "E39"; Value = "-3.45%" }
    @{ Cell = "D40"; Value = "0.04704" }
    @{ Cell = "E40"; Value = "-2.59%" }
    @{ Cell = "D41"; Value = "0.007518" }
    @{ Cell = "E41"; Value = "4.25%" }
    @{ Cell = "D42"; Value = "0.1364" }
    @{ Cell = "E42"; Value = "-0.75%" }
    @{ Cell = "D43"; Value = "0.002335" }
    @{ Cell = "E43"; Value = "9.98%" }
    @{ Cell = "D44"; Value = "0.01040" }
    @{ Cell = "E44"; Value = "4.60%" }
    @{ Cell = "D45"; Value = "0.00006168" }
    @{ Cell = "E45"; Value = "-1.62%" }
    @{ Cell = "D46"; Value = "0.00000000751" }
    @{ Cell = "E46"; Value = "0.06%" }
    @{ Cell = "D47"; Value = "1.552" }
    @{ Cell = "E47"; Value = "34.16%" }
    @{ Cell = "E48"; Value = "-16.75%" }
    @{ Cell = "D49"; Value = "0.00002102" }
    @{ Cell = "E49"; Value = "0.06%" }
    @{ Cell = "D50"; Value = "0.0002002" }
    @{ Cell = "E50"; Value = "0.06%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}

